# Update crypto price/volume data as scraped on Tue Jun  4 07:35:40 UTC 2024
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.972.67"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.18%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.769.21"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.20%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "630.52"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.00%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "165.50"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.43%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.767.30"
$ws.Range("D7").Style = "Normal"
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("E9").Value = "  +0.20%  "
$ws.Range("E10").Value = "  -1.54%  "
$ws.Range("E11").Value = "  +0.47%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.74"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.49%  "
$ws.Range("E13").Value = "  -4.01%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.81"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.97%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.405.16"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.19%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.770.54"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.02%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.971.19"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.13%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.67"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.84%  "
$ws.Range("E19").Value = "  +0.18%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.99"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.06%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "467.00"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.55%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.50"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.90%  "
$ws.Range("E23").Value = "  -0.99%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.14"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000143"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.63%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.07"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.41%  "
$ws.Range("E27").Value = "  -1.25%  "
$ws.Range("E28").Value = "  +0.72%  "
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.921.73"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.24%  "
$ws.Range("E31").Value = "  +0.57%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.26"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.84%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.09"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.34%  "
$ws.Range("E34").Value = "  +18.51%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "28.40"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.65%  "
$ws.Range("E36").Value = "  -0.09%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.723.51"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.05%  "
$ws.Range("E38").Value = "  -2.28%  "
$ws.Range("E39").Value = "  -0.81%  "
$ws.Range("E40").Value = "  -1.32%  "
$ws.Range("E41").Value = "  -2.33%  "
$ws.Range("E42").Value = "  +0.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.962"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.80%  "
$ws.Range("E45").Value = "  +5.29%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "156.27"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.02%  "
$ws.Range("B47").Value = "Arweave"
$ws.Range("C47").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "43.60"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.85%  "
$ws.Range("B48").Value = "ONDO"
$ws.Range("C48").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.41"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.78%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "46.92"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.21%  "
$ws.Range("E50").Value = "  -2.41%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.34"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.22%  "
